$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2068965517241379
$ws.Range("C2").Value = 0.5413793103448276
$ws.Range("J2").Value = 0.01379310344827586
$ws.Range("P2").Value = 0.1413793103448276
$ws.Range("S2").Value = 0.09655172413793103
$ws.Range("C3").Value = 0.01234567901234568
$ws.Range("J3").Value = 0.02469135802469136
$ws.Range("P3").Value = 0.7839506172839507
$ws.Range("S3").Value = 0.1790123456790123
$ws.Range("J4").Value = 0.07894736842105263
$ws.Range("P4").Value = 0.6578947368421053
$ws.Range("S4").Value = 0.2631578947368421
$ws.Range("B6").Value = 0.0411522633744856
$ws.Range("D6").Value = 0.00411522633744856
$ws.Range("F6").Value = 0.1111111111111111
$ws.Range("J6").Value = 0.2592592592592592
$ws.Range("O6").Value = 0.02469135802469136
$ws.Range("Q6").Value = 0.1481481481481481
$ws.Range("R6").Value = 0.06995884773662552
$ws.Range("S6").Value = 0.3415637860082305
$ws.Range("B7").Value = 0.1363636363636364
$ws.Range("D7").Value = 0.02272727272727273
$ws.Range("F7").Value = 0.07386363636363637
$ws.Range("J7").Value = 0.1022727272727273
$ws.Range("O7").Value = 0.005681818181818182
$ws.Range("Q7").Value = 0.1761363636363636
$ws.Range("R7").Value = 0.07954545454545454
$ws.Range("S7").Value = 0.4034090909090909
$ws.Range("B8").Value = 0.1005291005291005
$ws.Range("D8").Value = 0.01587301587301587
$ws.Range("F8").Value = 0.06349206349206349
$ws.Range("J8").Value = 0.08994708994708994
$ws.Range("O8").Value = 0.01851851851851852
$ws.Range("Q8").Value = 0.1111111111111111
$ws.Range("R8").Value = 0.1666666666666667
$ws.Range("S8").Value = 0.4338624338624338
$ws.Range("B9").Value = 0.0992063492063492
$ws.Range("D9").Value = 0.0119047619047619
$ws.Range("F9").Value = 0.1111111111111111
$ws.Range("J9").Value = 0.09126984126984126
$ws.Range("O9").Value = 0.007936507936507936
$ws.Range("Q9").Value = 0.1587301587301587
$ws.Range("R9").Value = 0.1388888888888889
$ws.Range("S9").Value = 0.3809523809523809
$ws.Range("B10").Value = 0.09969788519637462
$ws.Range("D10").Value = 0.01812688821752266
$ws.Range("E10").Value = 0.0007552870090634441
$ws.Range("F10").Value = 0.05891238670694864
$ws.Range("J10").Value = 0.1238670694864048
$ws.Range("O10").Value = 0.01963746223564955
$ws.Range("Q10").Value = 0.2009063444108761
$ws.Range("R10").Value = 0.1185800604229607
$ws.Range("S10").Value = 0.3595166163141994
$ws.Range("G11").Value = 0.1335740072202166
$ws.Range("J11").Value = 0.09025270758122744
$ws.Range("K11").Value = 0.1913357400722022
$ws.Range("L11").Value = 0.5667870036101083
$ws.Range("S11").Value = 0.01805054151624549
$ws.Range("G12").Value = 0.74375
$ws.Range("J12").Value = 0.1875
$ws.Range("K12").Value = 0.0125
$ws.Range("L12").Value = 0.0375
$ws.Range("S12").Value = 0.01875
$ws.Range("G13").Value = 0.6382978723404256
$ws.Range("J13").Value = 0.2978723404255319
$ws.Range("S13").Value = 0.06382978723404255
$ws.Range("F15").Value = 0.004739336492890996
$ws.Range("H15").Value = 0.1327014218009479
$ws.Range("I15").Value = 0.0947867298578199
$ws.Range("J15").Value = 0.3791469194312796
$ws.Range("K15").Value = 0.07582938388625593
$ws.Range("M15").Value = 0.01421800947867299
$ws.Range("O15").Value = 0.08530805687203792
$ws.Range("S15").Value = 0.2132701421800948
$ws.Range("F16").Value = 0.02150537634408602
$ws.Range("H16").Value = 0.1612903225806452
$ws.Range("I16").Value = 0.08064516129032258
$ws.Range("J16").Value = 0.4086021505376344
$ws.Range("K16").Value = 0.1236559139784946
$ws.Range("M16").Value = 0.02150537634408602
$ws.Range("N16").Value = 0.005376344086021506
$ws.Range("O16").Value = 0.05913978494623656
$ws.Range("S16").Value = 0.1182795698924731
$ws.Range("F17").Value = 0.01456310679611651
$ws.Range("H17").Value = 0.1286407766990291
$ws.Range("I17").Value = 0.133495145631068
$ws.Range("J17").Value = 0.424757281553398
$ws.Range("K17").Value = 0.09466019417475728
$ws.Range("M17").Value = 0.01941747572815534
$ws.Range("O17").Value = 0.05097087378640777
$ws.Range("S17").Value = 0.133495145631068
$ws.Range("F18").Value = 0.01048951048951049
$ws.Range("H18").Value = 0.1433566433566434
$ws.Range("I18").Value = 0.09090909090909091
$ws.Range("J18").Value = 0.4405594405594406
$ws.Range("K18").Value = 0.1048951048951049
$ws.Range("M18").Value = 0.02097902097902098
$ws.Range("O18").Value = 0.06643356643356643
$ws.Range("S18").Value = 0.1223776223776224
$ws.Range("F19").Value = 0.01803921568627451
$ws.Range("H19").Value = 0.1772549019607843
$ws.Range("I19").Value = 0.1058823529411765
$ws.Range("J19").Value = 0.3945098039215686
$ws.Range("K19").Value = 0.08705882352941176
$ws.Range("M19").Value = 0.0203921568627451
$ws.Range("N19").Value = 0.000784313725490196
$ws.Range("O19").Value = 0.06431372549019608
$ws.Range("S19").Value = 0.131764705882353
